$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.547.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.814.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.010'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.98%  '
$ws.Range("B5").Value = 'USDC'
$ws.Range("C5").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.007'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '308.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("E7").Value = '  -1.65%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3666'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07154'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8784'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07786'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.770.75'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.295'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("E15").Value = '  -2.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.14'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.16%  '
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008603'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.60%  '
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.634.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.14%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.012'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.984'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.54%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.084'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '112.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.862'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.83%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08685'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.059'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7348'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.486'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.119'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.02%  '
$ws.Range("B35").Value = 'Frax'
$ws.Range("C35").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.004'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.569'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.080'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01939'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05118'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.898'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.974'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5032'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1562'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.165'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.008'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4633'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.996'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.24%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '101.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.595'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06032'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.97%  '
